# Daily attendance processing - 2025-11-24 01:30:27
# Reorders the "Recorded By" email lists in column G for several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value  = "servinaz@med.asu.edu.eg, System, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
$ws.Range("G3").Value  = "System, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G4").Value  = "servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
$ws.Range("G5").Value  = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G7").Value  = "Kerelos.zareef@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, AbeerRagheb@med.asu.edu.eg"
$ws.Range("G15").Value = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"

$wb.Save()
